# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D (Price) and E (Volume 1h) keep their original text type
# (many values look numeric, e.g. "242.18", and Excel would otherwise silently
# convert them to numbers when assigned via .Value)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '36.277.72'
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").Value = '1.940.55'
$ws.Range("E3").Value = '  -4.19%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '242.18'
$ws.Range("E5").Value = '  -2.59%  '

$ws.Range("D6").Value = '0.607'
$ws.Range("E6").Value = '  -4.77%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").Value = '56.95'
$ws.Range("E8").Value = '  -9.23%  '

$ws.Range("D9").Value = '0.365'
$ws.Range("E9").Value = '  -7.03%  '

$ws.Range("D10").Value = '55.16'
$ws.Range("E10").Value = '  -4.97%  '

$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  +2.80%  '

$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  -0.25%  '

$ws.Range("D13").Value = '0.823'
$ws.Range("E13").Value = '  -8.28%  '

$ws.Range("D14").Value = '2.232.91'
$ws.Range("E14").Value = '  -3.82%  '

$ws.Range("D15").Value = '21.19'
$ws.Range("E15").Value = '  -9.63%  '

$ws.Range("D16").Value = '13.36'
$ws.Range("E16").Value = '  -7.05%  '

$ws.Range("D17").Value = '5.21'
$ws.Range("E17").Value = '  -6.07%  '

$ws.Range("D18").Value = '1.927.89'
$ws.Range("E18").Value = '  -4.64%  '

$ws.Range("D19").Value = '36.190.12'
$ws.Range("E19").Value = '  -1.35%  '

$ws.Range("D20").Value = '69.74'
$ws.Range("E20").Value = '  -3.38%  '

$ws.Range("D21").Value = '0.0₃0862'
$ws.Range("E21").Value = '  -2.64%  '

$ws.Range("D22").Value = '227.89'
$ws.Range("E22").Value = '  -3.68%  '

$ws.Range("D23").Value = '4.96'
$ws.Range("E23").Value = '  -7.88%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("D27").Value = '9.35'
$ws.Range("E27").Value = '  -5.46%  '

$ws.Range("D28").Value = '162.86'
$ws.Range("E28").Value = '  +1.81%  '

$ws.Range("D29").Value = '19.25'
$ws.Range("E29").Value = '  -5.87%  '

$ws.Range("E30").Value = '  -14.71%  '

$ws.Range("D31").Value = '0.117'
$ws.Range("E31").Value = '  -2.89%  '

$ws.Range("E32").Value = '  -3.40%  '

$ws.Range("D33").Value = '4.66'
$ws.Range("E33").Value = '  -8.15%  '

$ws.Range("D34").Value = '0.0626'
$ws.Range("E34").Value = '  +0.79%  '

$ws.Range("D35").Value = '4.27'
$ws.Range("E35").Value = '  -4.78%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = '6.02'
$ws.Range("E37").Value = '  -5.96%  '

$ws.Range("E38").Value = '  -1.93%  '

$ws.Range("D39").Value = '2.13'
$ws.Range("E39").Value = '  -11.15%  '

$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  -13.22%  '

$ws.Range("D41").Value = '0.0970'
$ws.Range("E41").Value = '  -2.96%  '

$ws.Range("E42").Value = '  -2.54%  '

$ws.Range("D43").Value = '1.17'
$ws.Range("E43").Value = '  -6.85%  '

$ws.Range("D44").Value = '0.0208'
$ws.Range("E44").Value = '  -3.40%  '

$ws.Range("D45").Value = '15.50'
$ws.Range("E45").Value = '  -9.03%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '7.34'
$ws.Range("E46").Value = '  -4.30%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.03'
$ws.Range("E47").Value = '  -9.52%  '

$ws.Range("D48").Value = '1.335.79'
$ws.Range("E48").Value = '  -2.73%  '

$ws.Range("D49").Value = '87.46'
$ws.Range("E49").Value = '  -6.89%  '

$ws.Range("D50").Value = '2.82'
$ws.Range("E50").Value = '  -2.95%  '

$ws.Range("D51").Value = '46.24'
$ws.Range("E51").Value = '  +2.13%  '

# Restore the default cell style so no stray number-format style is left behind
$ws.Range("B2:E51").Style = "Normal"